$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Set-ParagraphXml([string]$matchText, [string]$innerXml) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -eq $matchText) {
            $p.Range.InsertXML($pkgOpen + $innerXml + $pkgClose)
            return
        }
    }
    throw "paragraph not found: $matchText"
}

# 1) Sprint Number: 3 -> 2 (keep the three existing runs / rsids intact)
Set-ParagraphXml "Sprint Number: 3`r" '<w:p w14:paraId="70A84116" w14:textId="72A93904" w:rsidR="00F331F7" w:rsidRPr="00DD33C4" w:rsidRDefault="00F331F7" w:rsidP="00BB73E8"><w:pPr><w:rPr><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r><w:t>Sprint Number:</w:t></w:r><w:r w:rsidR="00BA54EE"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00931EC3"><w:t>2</w:t></w:r></w:p>'

# 2) "Writing unit tests" -> "Analyze the functions to write unit tests"
Set-ParagraphXml "Answer: Writing unit tests`r" '<w:p w14:paraId="0BDAF960" w14:textId="5C3B7EBB" w:rsidR="00F331F7" w:rsidRDefault="00F331F7" w:rsidP="00F331F7"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>Answer:</w:t></w:r><w:r w:rsidR="00BA54EE"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="0096450F"><w:t>Analyze the functions to write unit tests</w:t></w:r></w:p>'

# 3) "Continue" + " writing unit tests" -> single run "Analyze the functions to write unit tests"
Set-ParagraphXml "Answer: Continue writing unit tests`r" '<w:p w14:paraId="7F836323" w14:textId="459E295A" w:rsidR="00F331F7" w:rsidRDefault="00F331F7" w:rsidP="00F331F7"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>Answer:</w:t></w:r><w:r w:rsidR="00BA54EE"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="0096450F"><w:t>Analyze the functions to write unit tests</w:t></w:r></w:p>'

# 4) Styles.xml: remove <w:semiHidden/> from the DefaultParagraphFont character style
$styles = $d.Styles
$dpf = $styles.Item("Default Paragraph Font")
$dpf.Hidden = $false
